$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new 2021 data column (column R), mirroring the formatting of
# the equivalent cell in column Q (2020) for each row.
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("R4").Value = 2021

$ws.Range("Q5").Copy($ws.Range("R5"))
$ws.Range("R5").Value = 47.8

$ws.Range("Q6").Copy($ws.Range("R6"))
$ws.Range("R6").Value = 20.7

$ws.Range("Q7").Copy($ws.Range("R7"))
$ws.Range("R7").Value = 9.8

$ws.Range("Q8").Copy($ws.Range("R8"))
$ws.Range("R8").Value = 17.3

# Update the selected cell to reflect the new active selection in the sheet view
$ws.Range("P10").Select()
